$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Repulled/updated "dSF" (column F) data values for the relevant rows.
$updates = @{
    2  = 2
    3  = -2
    5  = -7
    6  = 0
    7  = -3
    10 = -5
    11 = -2
    12 = 1
    13 = 5
    14 = -1
    16 = -6
    17 = -1
    18 = -8
    19 = 11
    20 = -1
    21 = -4
    22 = -4
    24 = -3
    25 = -2
    26 = 2
    28 = -7
    29 = -1
    30 = -1
    31 = 2
    32 = -1
    33 = 5
    34 = -6
    35 = -2
    37 = -2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
